$p = $ppt.ActivePresentation

# --- Slide 2: add the missing space before the RTOS acronym expansion ---
# ("RTOS(Real Time Operating System)" -> "RTOS (Real Time Operating System)")
# NOTE: the run's apostrophe is a typographic right single quote (U+2019);
# the COM text bridge normalizes that to a plain "'" on read-back, so the
# replacement text is built explicitly with the correct code point instead
# of round-tripping through $run.Text.
$rsquo = [char]0x2019
$s2 = $p.Slides.Item(2)
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $shp = $s2.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.TextFrame.HasText) {
        $tr = $shp.TextFrame.TextRange
        $paraCount = $tr.Paragraphs().Count
        for ($j = 1; $j -le $paraCount; $j++) {
            $para = $tr.Paragraphs($j, 1)
            if ($para.Text -like "*RTOS(Real Time Operating System)*") {
                $runCount = $para.Runs().Count
                for ($k = 1; $k -le $runCount; $k++) {
                    $run = $para.Runs($k, 1)
                    if ($run.Text -like "*RTOS(Real Time Operating System)*") {
                        $run.Text = "These controller" + $rsquo + "s run RTOS (Real Time Operating System) "
                    }
                }
            }
        }
    }
}

# --- Remove the trailing "Intro to PLC Programming" slide ---
$last = $p.Slides.Item($p.Slides.Count)
$last.Delete()
